# Expand metadata tokenisation: rename "Assay Performer *" columns to
# "Assay Person *" on the isa_assay sheet, and restore the last-used
# selection (G11) on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_assay")

$ws.Range("A11").Value2 = "Assay Person Last Name"
$ws.Range("A12").Value2 = "Assay Person First Name"
$ws.Range("A13").Value2 = "Assay Person Mid Initials"
$ws.Range("A14").Value2 = "Assay Person Email"
$ws.Range("A15").Value2 = "Assay Person Phone"
$ws.Range("A16").Value2 = "Assay Person Fax"
$ws.Range("A17").Value2 = "Assay Person Address"
$ws.Range("A18").Value2 = "Assay Person Affiliation"
$ws.Range("A19").Value2 = "Assay Person Roles"
$ws.Range("A20").Value2 = "Assay Person Roles Term Accession Number"
$ws.Range("A21").Value2 = "Assay Person Roles Term Source REF"

# Reflect the last active selection recorded for this sheet in the saved file.
$ws.Activate()
[void]$ws.Range("G11").Select()
